$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin prices / volume figures as exact text values.
# Values are entered with a leading apostrophe to force Excel to keep
# them as literal text (preventing numeric auto-conversion / precision
# loss on values like "84.20" or "3.020"), then the style is reset back
# to "Normal" so no extra (quote-prefix) cell format lingers behind.

$ws.Range("D2").Value = "'28.953.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.95%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.905.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.28%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.43%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'324.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.31%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.38%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4597"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.58%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3814"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.35%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'45.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.69%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.26%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.9834"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.65%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -3.47%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.923.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.13%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.993"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.72%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.683"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.16%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.07064"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.77%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.39%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'84.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.04%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000009563"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.09%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'16.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.61%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.29%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'28.926.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.19%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.335"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.58%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'10.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.75%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.152.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.56%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.083"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.84%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'156.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.66%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.23%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.601"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -6.55%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'117.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.77%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.835"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.98%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09255"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.91%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.8599"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -5.30%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.107"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.71%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -7.10%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.020"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.99%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.05714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.98%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.148"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.02%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.004"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.02039"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.66%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.493"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.77%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.5529"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.94%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -3.73%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'9.354"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.74%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.734"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.47%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5201"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.11%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'11.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.55%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.095"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.12%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.06824"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.71%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'111.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.31%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'NEARProtocol"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'1.778"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.83%  "
$ws.Range("E51").Style = "Normal"
